# Insert a new price-check column at DQ (shifts existing "nom"/"url_produit"
# columns from DQ/DR to DR/DS), then fill the newly inserted column:
#  - header cell DQ1 gets the new check timestamp
#  - data rows (2..80) repeat the last known price from column DP (unchanged price)
#  - rows with no price on DP (81..206) are left blank, same as DP

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from column DQ onward one column to the right.
$ws.Columns("DQ").Insert()

# New header for the freshly inserted timestamp column.
$ws.Range("DQ1").Value() = "2026-02-02 07:43:04"

# Copy forward the last recorded price (column DP) into the new column DQ
# for every product row that currently has a price.
for ($r = 2; $r -le 80; $r++) {
    $price = $ws.Cells.Item($r, 120).Value()
    $ws.Cells.Item($r, 121).Value() = $price
}
